# "styled dataframe display mode"
# Replace the B-column literal values for rows 2-8 with their new
# (re-generated) figures, and append a brand-new row 9 ("Пастбище").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update existing column B values (rows 2-8) ---
$ws.Range("B2").Value = -25.55229949951172
$ws.Range("B3").Value = -70.23709869384766
$ws.Range("B4").Value = -56.36209869384766
$ws.Range("B5").Value = 90.31610107421875
$ws.Range("B6").Value = 93.73770141601562
$ws.Range("B7").Value = 22.84869956970215
$ws.Range("B8").Value = 30.17490005493164

# --- append new row 9 ("Пастбище") ---
$ws.Range("A9").Value = "Пастбище"
$ws.Range("B9").Value = 84.92591857910156
$ws.Range("C9").Value = 0.0108
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.354200005531311
$ws.Range("F9").Value = -0.07940000295639038
